# household_assessment.xlsx edit:
#  - survey sheet: insert a new "calculate" row (malaria_prone) right after
#    the existing visit_type calculate row (i.e. before the blank spacer
#    row that precedes "begin group household_assessment"). This pushes
#    every subsequent row down by one.
#  - fix a label typo on wire_mesh_on_inlets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Insert a new row at row 24 (shifts existing row 24.. down to 25..),
# inheriting formatting (style) from row 23 above it, same as Excel's
# default "Insert" behaviour for a selected row.
$ws.Range("A24").EntireRow.Insert()
$ws.Rows.Item(24).RowHeight = 15.75

$ws.Range("A24").Value2 = "calculate"
$ws.Range("B24").Value2 = "malaria_prone"
$ws.Range("C24").Value2 = "NO_LABEL"
$ws.Range("U24").Value2 = "../household_assessment/wire_mesh_on_inlets = 'yes' or ../household_assessment/stagnant_water = 'yes' or ../household_assessment/nets_availability = 'no' or ../household_assessment/proper_drainage = 'no' or ../household_assessment/resedual_spraying = 'no' or ../household_assessment/tall_grass = 'yes'"

# Fix label typo: "does have" -> "have" (row shifted to 27 after insert)
$ws.Range("C27").Value2 = "Does ${place_name}'s household have wire mesh on the windows/inlets?"
